$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 97, shifting old rows 97-99 down to 99-101
$ws.Rows.Item(97).Resize(2).Insert()

# Fill in new row 97 data
$ws.Cells.Item(97, 1).Value = 11
$ws.Cells.Item(97, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(97, 3).Value = "Bíobío"
$ws.Cells.Item(97, 4).Value = 44706
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 5).Value = 8
$ws.Cells.Item(97, 6).Value = 100112024
$ws.Cells.Item(97, 7).Value = "Choclo"
$ws.Cells.Item(97, 8).Value = "Choclero"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 4000
$ws.Cells.Item(97, 11).Value = 300
$ws.Cells.Item(97, 12).Value = 350
$ws.Cells.Item(97, 13).Value = 325
$ws.Cells.Item(97, 14).Value = "`$/unidad"
$ws.Cells.Item(97, 15).Value = "Región Metropolitana"
$ws.Cells.Item(97, 16).Value = 325
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# Fill in new row 98 data
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44706
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = 100112024
$ws.Cells.Item(98, 7).Value = "Choclo"
$ws.Cells.Item(98, 8).Value = "Choclero"
$ws.Cells.Item(98, 9).Value = "Segunda"
$ws.Cells.Item(98, 10).Value = 2000
$ws.Cells.Item(98, 11).Value = 250
$ws.Cells.Item(98, 12).Value = 250
$ws.Cells.Item(98, 13).Value = 250
$ws.Cells.Item(98, 14).Value = "`$/unidad"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 250
$ws.Cells.Item(98, 17).Value = 1
$ws.Cells.Item(98, 18).Value = "Hortaliza"
